# Update LR-pair TPM-derived values on Sheet1 per new TPM calculations.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 0.06515966666666667
$ws.Range("M2").Value = 4.717738333333333
$ws.Range("N2").Value = 14.153215
$ws.Range("O2").Value = 0.2002263444295212
$ws.Range("P2").Value = 0.2002263444295212
$ws.Range("Q2").Value = 0.3074062572205556
$ws.Range("R2").Value = 2.766656314985
$ws.Range("S2").Value = 0.2002263444295212
$ws.Range("T2").Value = 0.2002263444295212

# Row 3
$ws.Range("G3").Value = 0.06515966666666667
$ws.Range("N3").Value = 9.228847
$ws.Range("O3").Value = 0.1305610278731266
$ws.Range("P3").Value = 0.1305610278731266
$ws.Range("Q3").Value = 0.2004495314125556
$ws.Range("R3").Value = 1.804045782713
$ws.Range("S3").Value = 0.1305610278731266
$ws.Range("T3").Value = 0.1305610278731266

# Row 4
$ws.Range("G4").Value = 0.06515966666666667
$ws.Range("M4").Value = 0.6908423333333333
$ws.Range("N4").Value = 2.072527
$ws.Range("O4").Value = 0.02932015834857891
$ws.Range("P4").Value = 0.02932015834857891
$ws.Range("Q4").Value = 0.04501505615922222
$ws.Range("R4").Value = 0.405135505433
$ws.Range("S4").Value = 0.02932015834857891
$ws.Range("T4").Value = 0.02932015834857891

# Row 5
$ws.Range("G5").Value = 0.06515966666666667
$ws.Range("M5").Value = 15.077163
$ws.Range("N5").Value = 45.231489
$ws.Range("O5").Value = 0.6398924693487733
$ws.Range("P5").Value = 0.6398924693487733
$ws.Range("Q5").Value = 0.9824229153589999
$ws.Range("R5").Value = 8.841806238230999
$ws.Range("S5").Value = 0.6398924693487733
$ws.Range("T5").Value = 0.6398924693487733
